# Commit: "Add inheritance template & cfg formating"
#
# On the "Phiscal layout" sheet, insert a new "interface type" column
# between the "Device" column (A) and the "interface" column (old B),
# shifting interface/vlan/description one column to the right, and fill
# the new column with "fastethernet" for every existing device row.
#
# Cell-by-cell Copy (not Columns.Insert) is used so that only the columns
# that actually hold data are touched - Columns.Insert in this sheet would
# also manufacture brand-new empty cells on the blank placeholder rows
# (15-17), which the target layout does not have.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Phiscal layout")

# Rows that have 4 populated columns (A..D) -> need to become 5 (A..E).
# Shift existing data right one column, working from the rightmost column
# first so we never overwrite data before it has been copied onward.
$fullRows = 1,2,3,4,5,6,7,8,13
foreach ($r in $fullRows) {
    $ws.Range("D$r").Copy($ws.Range("E$r"))
    $ws.Range("C$r").Copy($ws.Range("D$r"))
    $ws.Range("B$r").Copy($ws.Range("C$r"))
}

# Rows that only have 2 populated columns (A..B) -> need to become 3 (A..C).
$shortRows = 9,10,11,12
foreach ($r in $shortRows) {
    $ws.Range("B$r").Copy($ws.Range("C$r"))
}

# Header for the newly inserted column.
$ws.Range("B1").Value = "interface type"

# Data rows all use the same interface type.
$dataRows = 2,3,4,5,6,7,8,9,10,11,12,13
foreach ($r in $dataRows) {
    $ws.Range("B$r").Value = "fastethernet"
}

# Rows 15-17 are blank placeholder rows that only had a styled, empty A/B
# pair. Move the empty, styled cell from B to C and drop B entirely (no
# new column-type cell is introduced on these blank rows).
$blankRows = 15,16,17
foreach ($r in $blankRows) {
    $ws.Range("B$r").Copy($ws.Range("C$r"))
    $ws.Range("B$r").Clear()
}
